# Split runs that contain inline <exp>...</exp> (and <corr>...</corr>)
# pseudo-tags out of their surrounding plain text into separate Word
# runs, so the tag markers get their own (small, grey, Courier New)
# character formatting - mirroring how the rest of the document already
# renders its <tag> markers.
#
# For every edit we:
#   1. Find.Execute the *original* (still-merged) text to anchor a Range.
#   2. Re-express that whole span as $d.Range(start, end) sub-ranges, one
#      per literal slice (tag text vs. plain text).
#   3. Apply Font.Name / Font.Size / Font.Color to the tag slices only -
#      Word automatically splits the underlying run(s) for us whenever a
#      sub-range's formatting diverges from its neighbours.

$d = $word.ActiveDocument

function Format-TagRun($range) {
    $range.Font.Name = "Courier New"
    $range.Font.Size = 7
    $range.Font.Color = 11119017   # a9a9a9
}

function Format-CorrRun($range) {
    $range.Font.Name = "Courier New"
    $range.Font.Size = 9
    $range.Font.Color = 1118633    # a91111
}

function Split-ExpRun($searchText, $pieces) {
    # $pieces is an array of @{ Text = "..."; Tag = $true/$false }
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false,
                                $false, $false, $true, 1, $false,
                                "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return
    }
    $s = $rng.Start
    $pos = $s
    foreach ($piece in $pieces) {
        $len = $piece.Text.Length
        $sub = $d.Range($pos, $pos + $len)
        if ($piece.Tag -eq "tag") {
            Format-TagRun $sub
        } elseif ($piece.Tag -eq "corr") {
            Format-CorrRun $sub
        }
        $pos = $pos + $len
    }
}

# 1) co<exp>mm</exp>e dict est de
Split-ExpRun "co<exp>mm</exp>e dict est de " @(
    @{ Text = "co"; Tag = "plain" },
    @{ Text = "<exp>"; Tag = "tag" },
    @{ Text = "mm"; Tag = "plain" },
    @{ Text = "</exp>"; Tag = "tag" },
    @{ Text = "e dict est de "; Tag = "plain" }
)

# 2) en traina<exp>n</exp>t
Split-ExpRun "en traina<exp>n</exp>t" @(
    @{ Text = "en traina"; Tag = "plain" },
    @{ Text = "<exp>"; Tag = "tag" },
    @{ Text = "n"; Tag = "plain" },
    @{ Text = "</exp>"; Tag = "tag" },
    @{ Text = "t"; Tag = "plain" }
)

# 3) co<exp>mm</exp>e deulx
Split-ExpRun "co<exp>mm</exp>e deulx " @(
    @{ Text = "co"; Tag = "plain" },
    @{ Text = "<exp>"; Tag = "tag" },
    @{ Text = "mm"; Tag = "plain" },
    @{ Text = "</exp>"; Tag = "tag" },
    @{ Text = "e deulx "; Tag = "plain" }
)

# 4)  esta<exp>n</exp>t
Split-ExpRun " esta<exp>n</exp>t" @(
    @{ Text = " esta"; Tag = "plain" },
    @{ Text = "<exp>"; Tag = "tag" },
    @{ Text = "n"; Tag = "plain" },
    @{ Text = "</exp>"; Tag = "tag" },
    @{ Text = "t"; Tag = "plain" }
)

# 5) ...doulcem<exp>ent</exp> et si par cas
Split-ExpRun "hault pource quen ceste sorte ilz se deseichent doulcem<exp>ent</exp> et si par cas" @(
    @{ Text = "hault pource quen ceste sorte ilz se deseichent doulcem"; Tag = "plain" },
    @{ Text = "<exp>"; Tag = "tag" },
    @{ Text = "ent"; Tag = "plain" },
    @{ Text = "</exp>"; Tag = "tag" },
    @{ Text = " et si par cas"; Tag = "plain" }
)

# 6) ...doulcem<corr><exp>ent</exp></corr>
Split-ExpRun "cest quand il est fondu doulcem<corr><exp>ent</exp></corr>" @(
    @{ Text = "cest quand il est fondu doulcem"; Tag = "plain" },
    @{ Text = "<corr>"; Tag = "corr" },
    @{ Text = "<exp>"; Tag = "tag" },
    @{ Text = "ent"; Tag = "plain" },
    @{ Text = "</exp>"; Tag = "tag" },
    @{ Text = "</corr>"; Tag = "corr" }
)
